$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text (Coinranking price strings
# such as "31.80" or "1.00" look numeric to Excel and would otherwise be
# auto-converted, dropping the trailing zero). Mark the range as Text first.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

$ws.Range("D2").Value = "34.578.53"
$ws.Range("E2").Value = "  +13.20%  "

$ws.Range("D3").Value = "1.827.19"
$ws.Range("E3").Value = "  +9.25%  "

Set-TextValue $ws.Range("D4") "0.997"
$ws.Range("E4").Value = "  -0.23%  "

Set-TextValue $ws.Range("D5") "230.79"
$ws.Range("E5").Value = "  +5.14%  "

$ws.Range("E6").Value = "  +4.27%  "

Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.15%  "

Set-TextValue $ws.Range("D8") "31.80"
$ws.Range("E8").Value = "  +7.63%  "

Set-TextValue $ws.Range("D9") "47.11"
$ws.Range("E9").Value = "  +6.91%  "

$ws.Range("E10").Value = "  +7.22%  "

$ws.Range("E11").Value = "  +5.11%  "

Set-TextValue $ws.Range("D12") "0.0932"
$ws.Range("E12").Value = "  +2.95%  "

$ws.Range("D13").Value = "2.087.92"
$ws.Range("E13").Value = "  +9.17%  "

$ws.Range("D14").Value = "1.828.75"
$ws.Range("E14").Value = "  +9.26%  "

Set-TextValue $ws.Range("D15") "0.654"
$ws.Range("E15").Value = "  +6.53%  "

$ws.Range("D16").Value = "34.448.14"
$ws.Range("E16").Value = "  +12.68%  "

Set-TextValue $ws.Range("D17") "10.43"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("E18").Value = "  +7.26%  "

Set-TextValue $ws.Range("D19") "69.74"
$ws.Range("E19").Value = "  +5.09%  "

Set-TextValue $ws.Range("D20") "259.30"
$ws.Range("E20").Value = "  +6.74%  "

$ws.Range("D21").Value = "0.0₃0755"
$ws.Range("E21").Value = "  +4.53%  "

$ws.Range("E22").Value = "  -0.07%  "

Set-TextValue $ws.Range("D23") "10.58"
$ws.Range("E23").Value = "  +5.81%  "

Set-TextValue $ws.Range("D24") "4.37"
$ws.Range("E24").Value = "  +2.04%  "

$ws.Range("E25").Value = "  +2.81%  "

Set-TextValue $ws.Range("D26") "158.95"
$ws.Range("E26").Value = "  +0.35%  "

Set-TextValue $ws.Range("D27") "16.71"
$ws.Range("E27").Value = "  +5.43%  "

$ws.Range("E28").Value = "  +6.99%  "

$ws.Range("E29").Value = "  +2.45%  "

Set-TextValue $ws.Range("D30") "1.00"
$ws.Range("E30").Value = "  +0.15%  "

Set-TextValue $ws.Range("D31") "3.89"
$ws.Range("E31").Value = "  +12.31%  "

Set-TextValue $ws.Range("D32") "0.0521"
$ws.Range("E32").Value = "  +5.16%  "

Set-TextValue $ws.Range("D33") "1.22"
$ws.Range("E33").Value = "  +6.61%  "

Set-TextValue $ws.Range("D34") "3.59"
$ws.Range("E34").Value = "  +9.12%  "

$ws.Range("D35").Value = "1.553.83"
$ws.Range("E35").Value = "  +4.34%  "

Set-TextValue $ws.Range("D36") "1.80"
$ws.Range("E36").Value = "  +2.02%  "

$ws.Range("E37").Value = "  +6.27%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D38") "0.637"
$ws.Range("E38").Value = "  +6.65%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D39") "0.0190"
$ws.Range("E39").Value = "  +6.91%  "

$ws.Range("B40").Value = "MinaProtocolToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
Set-TextValue $ws.Range("D40") "1.27"
$ws.Range("E40").Value = "  +208.52%  "

Set-TextValue $ws.Range("D41") "84.94"
$ws.Range("E41").Value = "  +0.65%  "

Set-TextValue $ws.Range("D42") "2.81"
$ws.Range("E42").Value = "  +5.21%  "

Set-TextValue $ws.Range("D43") "0.920"
$ws.Range("E43").Value = "  +9.64%  "

Set-TextValue $ws.Range("D44") "2.32"
$ws.Range("E44").Value = "  +1.38%  "

$ws.Range("E45").Value = "  +10.44%  "

$ws.Range("E46").Value = "  +5.53%  "

$ws.Range("E47").Value = "  +5.46%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D48") "12.58"
$ws.Range("E48").Value = "  +29.11%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.990.15"
$ws.Range("E49").Value = "  +10.09%  "

$ws.Range("E50").Value = "  +5.42%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D51") "53.21"
$ws.Range("E51").Value = "  +4.23%  "
